$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 474; everything currently at/after row 474
# shifts down by one (old row 474 -> 475, ..., old row 498 -> 499).
$ws.Rows(474).Insert()

# Populate the newly inserted row with the new weekly price record.
$ws.Cells.Item(474, 1).Value = 6
$ws.Cells.Item(474, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(474, 3).Value = "Metropolitana"
$ws.Cells.Item(474, 4).Value = 44753
$ws.Cells.Item(474, 5).Value = 13
$ws.Cells.Item(474, 6).Value = 100112052
$ws.Cells.Item(474, 7).Value = "Albahaca"
$ws.Cells.Item(474, 8).Value = "Sin especificar"
$ws.Cells.Item(474, 9).Value = "Primera"
$ws.Cells.Item(474, 10).Value = 80
$ws.Cells.Item(474, 11).Value = 3500
$ws.Cells.Item(474, 12).Value = 4000
$ws.Cells.Item(474, 13).Value = 3812
$ws.Cells.Item(474, 14).Value = "`$/paquete"
$ws.Cells.Item(474, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(474, 16).Value = 3812
$ws.Cells.Item(474, 17).Value = 1
$ws.Cells.Item(474, 18).Value = "Hortaliza"
